$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": drop the two oldest weekly rows -----------
# The original rows 2 ("2023-07-16") and 3 ("2023-07-23") are removed; the
# remaining weekly rows shift up, so the sheet shrinks from A1:B7 to A1:B5.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(2).Delete()
$wsWeekly.Rows.Item(2).Delete()

# --- Sheet "Monthly Trend": July total drops from 10 to 8 ---------------
# Reflects the two removed weekly rows (quantities 1 + 1 = 2 less).
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(2, 2).Value = 8
